# Updated capital structure database
# Apply refreshed metrics to the Belgium "Software (Entertainment)" rows
# (row 2: aggregate row, row 3: AudioValley SA). Both data rows receive
# the same refreshed values per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G"  = -0.1136170212765958
    "H"  = -0.1782978723404255
    "I"  = -0.2182978723404255
    "J"  = -0.2088304735758408
    "K"  = 5.65
    "L"  = 0.2404255319148936
    "M"  = 0.008
    "N"  = 0.0001724137931034483
    "O"  = 0.001415929203539823
    "R"  = 0
    "S"  = 0.008
    "U"  = 0.833
    "V"  = 0.01795258620689655
    "W"  = 0.8814352574102965
    "X"  = 0.08815409990444158
    "Y"  = 0.7932811575058549
    "Z"  = 1.140721324207563
    "AA" = -0.2382173743523255
    "AB" = 0.06640739664943104
    "AC" = -0.3046247710017566
    "AD" = 20.6
    "AE" = 0
    "AF" = 20.6
    "AG" = 19.767
    "AH" = 0.3074626865671642
    "AI" = 0.5309278350515465
    "AJ" = 0.298744086931552
    "AK" = 0.520636342086549
    "AL" = 2.07
    "AM" = 2.034
    "AN" = -3.570190641247834
    "AO" = -2.478260869565218
    "AP" = -3.425823223570192
    "AQ" = -2.52212389380531
}

foreach ($row in 2,3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
